$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every value in the Price/Volume columns as
# literal text (including numeric- and percent-looking strings). A plain
# Range.Value assignment would let Excel "smart parse" these into real
# numbers/percentages, so each target cell is first (re)marked as Text
# (`@`) before its literal string is written.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "262.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.90%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.71%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.678"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06106"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.05%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.708"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.69%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8496"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.25%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9094"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.87%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.04907"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.15%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07091"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.10%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03113"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.03%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09042"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.06%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001524"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.09%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0006192"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.44%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005984"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.28%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.01%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.169"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.46%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.146"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.54%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1409"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.21%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.52%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.081"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.21%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04228"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.03%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001176"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.32%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004060"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.78%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.03%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "23.08%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03929"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.45%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.13%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004175"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.08%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.27%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01161"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-28.90%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005095"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.11%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.02%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2582"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "109.99%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.02%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.02%"
